$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organizers")

# Fill in the "imageName" column (I) for each organizer row with their photo filename.
$ws.Range("I2").Value  = "camden_squire.jpg"
$ws.Range("I3").Value  = "guangze_zu.jpg"
$ws.Range("I4").Value  = "david_richey.jpg"
$ws.Range("I5").Value  = "gautam_sapre.jpg"
$ws.Range("I6").Value  = "pallavi_vayalali.jpg"
$ws.Range("I7").Value  = "caitlin_tibbetts.jpg"
$ws.Range("I8").Value  = "sivam_patel.jpg"
$ws.Range("I9").Value  = "soham_mukerjee.jpg"
$ws.Range("I10").Value = "chaithu_dikkala.jpg"
$ws.Range("I11").Value = "neha_rode.jpg"
$ws.Range("I12").Value = "aditya_guin.jpg"
$ws.Range("I13").Value = "austin_luong.jpg"
$ws.Range("I14").Value = "rolando_gonzalez.jpg"
$ws.Range("I15").Value = "atharv_jain.jpg"
$ws.Range("I16").Value = "willie_chalmers iii.jpg"
$ws.Range("I17").Value = "sanjana_sivakumar.jpg"
$ws.Range("I18").Value = "vishvak_bandi.jpg"
$ws.Range("I19").Value = "alexander_osypov.jpg"
$ws.Range("I20").Value = "michael_xu.jpg"
$ws.Range("I21").Value = "daniel_wang.jpg"
$ws.Range("I22").Value = "ishpreet_bhasin.jpg"
$ws.Range("I23").Value = "elvina_almeida.jpg"
$ws.Range("I24").Value = "abhitej_arora.jpg"

# I16 previously carried a stray hyperlink-style font, and I17:I24 are brand
# new cells with no style at all; normalize all of them in one pass to the
# same plain black Arial 10 look used by the rest of the imageName column
# (I2:I15).
$ws.Range("I16:I24").Font.ThemeColor = 1
$ws.Range("I16:I24").Font.Underline = $false

# Move the active selection to I3, matching where the author left off editing.
$ws.Range("I3").Select() | Out-Null
